$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date text (keep as text, same format as existing date-like strings in the sheet)
$ws.Range("N2").Value = "2018-12-31 00:00:00"

# Numeric financial figures
$ws.Range("O2").Value = 529139010.02
$ws.Range("P2").Value = 118485052.82
$ws.Range("Q2").Value = 11054829.8
$ws.Range("R2").Value = 7.4436036534
$ws.Range("S2").Value = 127964434.06
$ws.Range("T2").Value = 124.5551178132
$ws.Range("U2").Value = 139602091.71
$ws.Range("V2").Value = 47.3383802008
$ws.Range("W2").Value = 262909016.93
$ws.Range("X2").Value = 139907186.02
$ws.Range("Y2").Value = 53.7184450887
$ws.Range("Z2").Value = 7173597.07
$ws.Range("AA2").Value = -17.448418909
$ws.Range("AB2").Value = 266229993.09
$ws.Range("AC2").Value = 20.591176787
$ws.Range("AD2").Value = 34.3684446973
$ws.Range("AE2").Value = 51.947367108
$ws.Range("AF2").Value = 129.4633296613
$ws.Range("AG2").Value = 49.6861905759
